# Guns n Roses - Sweet Child of mine
# "Fix back image bug": the two lyric slides get collapsed down to short
# placeholder captions, and a third slide (a copy of the old slide 2,
# also collapsed to a placeholder caption) is appended at the end.

$p = $ppt.ActivePresentation

# --- Step 1: duplicate slide 2 (still holding its full lyrics) so the
# new slide 3 starts life as a faithful copy of slide 2, exactly like
# the new slide shows up in the authored diff. This must happen before
# slide 2's own text gets trimmed down below. ---
$s2 = $p.Slides.Item(2)
$dupRange = $s2.Duplicate()
$s3 = $dupRange.Item(1)

# --- Step 2: trim slide 1's lyrics down to a short caption ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = "GNR Sweet Child of mine 1`r"

# --- Step 3: trim slide 2's lyrics down to a short caption, bumping the
# font size from 16pt to 20pt like the rest of the shape text. ---
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "GNR Sweet Child of mine 2`r"
$tr2.Font.Size = 20

# --- Step 4: trim slide 3's (the new slide's) lyrics down to a short
# caption split across two runs, also at 20pt. ---
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "GNR Sweet Child of mine 3`r"
$tr3.Font.Size = 20
$run3b = $tr3.Characters(20, 6)
$run3b.Font.Bold = $true
